$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "median"
$ws.Range("B3").Value = "AstorMain-median-3"
$ws.Range("C3").Value = "/Users/ruizhengu/Experiments/APR-as-AAT/IntroClassJava/dataset/median/fe9d5fb933c758c2cbbd859e3ecbd2d2eaecc2843b57cfd97da99af24de59f189a144d13ce81ec585d7c2f7367f70c0fb2aec8269eed1fd8380def614817ef7c/002"
